$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 16
$ws.Range("C7").Value = 20
$ws.Range("C8").Value = 14
$ws.Range("C9").Value = 12
$ws.Range("C10").Value = 13
$ws.Range("C12").Value = 10
$ws.Range("C15").Value = 12
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 14
$ws.Range("C18").Value = 15
